# explanation.pptx - "Fixed even more typos (in the PPT)"
#
# 1. Slide 2  : Title "Software specitication" -> "Software specification"
# 2. Slide 5  : Title "pp03.ppy: ..."          -> "pp03.py: ..."
# 3. Slide 6  : Title "pp04.ppy: ..."          -> "pp04.py: ..."
# 4. Slide 7  : Title "pp05.ppy: ..."          -> "pp05.py: ..."
# 5. Slide 9  : removed entirely (the "Example of Flask in action" / pitdb.org slide)

$p = $ppt.ActivePresentation

# --- Slide 2: title run is split across "Software " / "specitication" runs,
# so delete the whole range and retype it to collapse it back into a single,
# clean run (matches the author's corrected single-run title). ---
$s2 = $p.Slides.Item(2)
$titleRange2 = $s2.Shapes.Item(1).TextFrame.TextRange
$titleRange2.Delete()
$titleRange2.InsertAfter("Software specification") | Out-Null

# --- Slides 5-7: titles are already single runs, so a straight text
# assignment fixes the "ppy" -> "py" typo without disturbing formatting. ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "pp03.py: Return some useful data using HTML"

$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "pp04.py: Use templates for cleaner code"

$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "pp05.py: Accept user input by using a form "

# --- Slide 9 ("Example of Flask in action" / pitdb.org) is deleted outright. ---
$p.Slides.Item(9).Delete() | Out-Null
